$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.661.37'
$ws.Range("E2").Value = '  +1.16%  '

$ws.Range("D3").Value = '3.464.72'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '414.99'
$ws.Range("E5").Value = '  +1.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.33'
$ws.Range("E6").Value = '  +1.48%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  -0.74%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.730'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.140'
$ws.Range("E10").Value = '  +0.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.88'
$ws.Range("E11").Value = '  +0.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.73'
$ws.Range("E12").Value = '  +6.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000218'
$ws.Range("E13").Value = '  +5.15%  '

$ws.Range("D14").Value = '4.013.60'
$ws.Range("E14").Value = '  +1.16%  '

$ws.Range("E15").Value = '  -0.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.54'
$ws.Range("E16").Value = '  -4.32%  '

$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.83'
$ws.Range("E17").Value = '  +3.52%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.446.46'
$ws.Range("E18").Value = '  +1.10%  '

$ws.Range("E19").Value = '  -0.15%  '

$ws.Range("D20").Value = '62.560.92'
$ws.Range("E20").Value = '  +0.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '470.47'
$ws.Range("E21").Value = '  +5.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '90.86'
$ws.Range("E22").Value = '  -0.50%  '

$ws.Range("E23").Value = '  +2.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.45'
$ws.Range("E24").Value = '  +3.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.51'
$ws.Range("E25").Value = '  +20.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.32'
$ws.Range("E26").Value = '  +2.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.37'
$ws.Range("E27").Value = '  +1.78%  '

$ws.Range("E28").Value = '  +0.75%  '

$ws.Range("E29").Value = '  -2.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.05'
$ws.Range("E30").Value = '  +0.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.75'
$ws.Range("E31").Value = '  +1.03%  '

$ws.Range("E32").Value = '  -2.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.112'
$ws.Range("E33").Value = '  -1.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.19'
$ws.Range("E34").Value = '  -4.08%  '

$ws.Range("E35").Value = '  +0.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.67'
$ws.Range("E36").Value = '  +9.61%  '

$ws.Range("E37").Value = '  -1.71%  '

$ws.Range("E38").Value = '  +0.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.06'
$ws.Range("E39").Value = '  +4.38%  '

$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("E41").Value = '  -0.45%  '

$ws.Range("E42").Value = '  +0.28%  '

$ws.Range("E43").Value = '  +6.70%  '

$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '144.59'
$ws.Range("E44").Value = '  +1.80%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.39'
$ws.Range("E45").Value = '  +3.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.07'
$ws.Range("E46").Value = '  +4.56%  '

$ws.Range("E47").Value = '  +12.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.40'
$ws.Range("E49").Value = '  -1.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.25'
$ws.Range("E50").Value = '  -0.50%  '

$ws.Range("E51").Value = '  -2.14%  '
